$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Keywords"

# Rename the table (this also updates structured-reference formulas automatically)
$lo = $ws.ListObjects.Item(1)
$lo.Name = "Cluster_Keywords"

for ($r = 2; $r -le 28; $r++) {
    $ws.Range("B$r").Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
}

$wb.Save()
